$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 80019
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = $null
# Row 23
$ws.Range("H23").Value = 80019
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = $null
# Row 29
$ws.Range("H29").Value = 799.6667
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = $null
# Row 38
$ws.Range("H38").Value = 3391.8
$ws.Range("I38").Value = 84.625
$ws.Range("J38").Value = 7171.4287
$ws.Range("K38").Value = 253.875
$ws.Range("L38").Value = 21514.2861
$ws.Range("M38").Value = 118.125
$ws.Range("N38").Value = -22258.2861
# Row 54
$ws.Range("H54").Value = 41999.715
# Row 58
$ws.Range("H58").Value = 5968.0557
$ws.Range("I58").Value = 552.5
$ws.Range("J58").Value = 12737.5
$ws.Range("K58").Value = 1657.5
$ws.Range("L58").Value = 38212.5
$ws.Range("M58").Value = -1507.5
$ws.Range("N58").Value = -38512.5
# Row 76
$ws.Range("H76").Value = 3156.4783
$ws.Range("I76").Value = 3139.95
$ws.Range("J76").Value = 3266.6667
$ws.Range("K76").Value = 3139.95
$ws.Range("L76").Value = 3266.6667
$ws.Range("M76").Value = -2824.95
$ws.Range("N76").Value = -3896.6667
# Row 79
$ws.Range("H79").Value = 3156.4783
$ws.Range("I79").Value = 3139.95
$ws.Range("J79").Value = 3266.6667
$ws.Range("K79").Value = 3139.95
$ws.Range("L79").Value = 3266.6667
$ws.Range("M79").Value = -2047.95
$ws.Range("N79").Value = -5450.6667
# Row 86
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null
# Row 89
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null
# Row 129
$ws.Range("H129").Value = 944.3728599999999
$ws.Range("J129").Value = 979.1786
$ws.Range("L129").Value = 2937.5358
$ws.Range("N129").Value = -12937.5358
# Row 137
$ws.Range("H137").Value = 3440
$ws.Range("I137").Value = 2364.4443
$ws.Range("J137").Value = 6666.6665
$ws.Range("K137").Value = 7093.3329
$ws.Range("L137").Value = 19999.9995
$ws.Range("M137").Value = -4543.3329
$ws.Range("N137").Value = -25099.9995
# Row 138
$ws.Range("H138").Value = 3163.1943
$ws.Range("I138").Value = 1856.1666
$ws.Range("J138").Value = 3424.6
$ws.Range("K138").Value = 5568.4998
$ws.Range("L138").Value = 10273.8
$ws.Range("M138").Value = -428.4997999999996
$ws.Range("N138").Value = -20553.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5951.154
$ws.Range("I32").Value = 3942.558
$ws.Range("K32").Value = 3942.558
$ws.Range("M32").Value = -3655.558
# Row 61
$ws.Range("H61").Value = 2270.4783
$ws.Range("I61").Value = 1799.909
$ws.Range("J61").Value = 2701.8333
$ws.Range("K61").Value = 1799.909
$ws.Range("L61").Value = 2701.8333
$ws.Range("M61").Value = -1587.909
$ws.Range("N61").Value = -3125.8333
# Row 136
$ws.Range("H136").Value = 2270.4783
$ws.Range("I136").Value = 1799.909
$ws.Range("J136").Value = 2701.8333
$ws.Range("K136").Value = 5399.727000000001
$ws.Range("L136").Value = 8105.499899999999
$ws.Range("M136").Value = -2849.727000000001
$ws.Range("N136").Value = -13205.4999
# Row 137
$ws.Range("H137").Value = 39586.668
$ws.Range("J137").Value = 39586.668
$ws.Range("L137").Value = 39586.668
$ws.Range("N137").Value = -49786.668

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 7131.6924
$ws.Range("I20").Value = 1400.3334
$ws.Range("J20").Value = 20027.25
$ws.Range("K20").Value = 1400.3334
$ws.Range("L20").Value = 20027.25
$ws.Range("M20").Value = -1153.3334
$ws.Range("N20").Value = -20521.25
# Row 59
$ws.Range("H59").Value = 64816.332
$ws.Range("J59").Value = 64816.332
$ws.Range("L59").Value = 64816.332
$ws.Range("N59").Value = -66510.33199999999
# Row 105
$ws.Range("H105").Value = 2429.121
$ws.Range("I105").Value = 2354.862
$ws.Range("J105").Value = 2967.5
$ws.Range("K105").Value = 2354.862
$ws.Range("L105").Value = 2967.5
$ws.Range("M105").Value = -607.8620000000001
$ws.Range("N105").Value = -6461.5
# Row 137
$ws.Range("H137").Value = 25000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null
# Row 129
$ws.Range("H129").Value = 5882.5
$ws.Range("I129").Value = 5882.5
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 17647.5
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -12647.5
$ws.Range("N129").Value = $null
# Row 136
$ws.Range("H136").Value = 3520
$ws.Range("I136").Value = 3306.6667
$ws.Range("J136").Value = 4800
$ws.Range("K136").Value = 9920.000100000001
$ws.Range("L136").Value = 14400
$ws.Range("M136").Value = -4820.000100000001
$ws.Range("N136").Value = -24600
# Row 137
$ws.Range("H137").Value = 3745.95
$ws.Range("J137").Value = 4507.4375
$ws.Range("L137").Value = 13522.3125
$ws.Range("N137").Value = -23722.3125
# Row 138
$ws.Range("H138").Value = 3139.5
$ws.Range("I138").Value = 2249.1667
$ws.Range("J138").Value = 4475
$ws.Range("K138").Value = 6747.500100000001
$ws.Range("L138").Value = 13425
$ws.Range("M138").Value = -1607.500100000001
$ws.Range("N138").Value = -23705
# Row 139
$ws.Range("H139").Value = 1384
$ws.Range("I139").Value = 1204.4445
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 3613.3335
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 1526.6665
$ws.Range("N139").Value = -19280
# Row 140
$ws.Range("H140").Value = 19943.31
$ws.Range("I140").Value = 28861.895
$ws.Range("J140").Value = 2998
$ws.Range("K140").Value = 86585.685
$ws.Range("L140").Value = 8994
$ws.Range("M140").Value = -81405.685
$ws.Range("N140").Value = -19354

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 30091
$ws.Range("J46").Value = 30091
$ws.Range("L46").Value = 30091
$ws.Range("N46").Value = -30403
# Row 70
$ws.Range("H70").Value = 5602.268
$ws.Range("I70").Value = 5015.816
$ws.Range("K70").Value = 5015.816
$ws.Range("M70").Value = -4745.816
# Row 73
$ws.Range("H73").Value = 5602.268
$ws.Range("I73").Value = 5015.816
$ws.Range("K73").Value = 5015.816
$ws.Range("M73").Value = -4079.816
# Row 97
$ws.Range("H97").Value = 1554.8334
$ws.Range("I97").Value = 1332.25
$ws.Range("K97").Value = 1332.25
$ws.Range("M97").Value = -836.25
# Row 102
$ws.Range("H102").Value = 1786.9656
$ws.Range("J102").Value = 2605.5
$ws.Range("L102").Value = 2605.5
$ws.Range("N102").Value = -5849.5
# Row 122
$ws.Range("H122").Value = 4772.273
$ws.Range("I122").Value = 2537.4
$ws.Range("K122").Value = 7612.200000000001
$ws.Range("M122").Value = -5162.200000000001
# Row 123
$ws.Range("H123").Value = 11569.2
$ws.Range("J123").Value = 11569.2
$ws.Range("L123").Value = 11569.2
$ws.Range("N123").Value = -16469.2
# Row 132
$ws.Range("H132").Value = 3267.2812
$ws.Range("I132").Value = 2397.8667
$ws.Range("K132").Value = 7193.6001
$ws.Range("M132").Value = -4663.6001
# Row 137
$ws.Range("H137").Value = 42658.332
$ws.Range("J137").Value = 42658.332
$ws.Range("L137").Value = 42658.332
$ws.Range("N137").Value = -52858.332

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = $null
# Row 46
$ws.Range("H46").Value = 3540.4
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 2567.3333
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 2567.3333
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -2943.3333
# Row 50
$ws.Range("H50").Value = 37359
$ws.Range("J50").Value = 37359
$ws.Range("L50").Value = 37359
$ws.Range("N50").Value = -38633
# Row 54
$ws.Range("H54").Value = 35082.5
$ws.Range("J54").Value = 35082.5
$ws.Range("L54").Value = 35082.5
$ws.Range("N54").Value = -36370.5
# Row 136
$ws.Range("H136").Value = 3525.6155
$ws.Range("I136").Value = 1666.6316
$ws.Range("K136").Value = 4999.8948
$ws.Range("M136").Value = -2449.8948

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 9450
$ws.Range("I45").Value = 5000
$ws.Range("J45").Value = 10933.333
$ws.Range("K45").Value = 5000
$ws.Range("L45").Value = 10933.333
$ws.Range("M45").Value = -4509
$ws.Range("N45").Value = -11915.333
# Row 122
$ws.Range("H122").Value = 3781.6667
$ws.Range("I122").Value = 2082.3076
$ws.Range("J122").Value = 8200
$ws.Range("K122").Value = 6246.9228
$ws.Range("L122").Value = 24600
$ws.Range("M122").Value = -3796.9228
$ws.Range("N122").Value = -29500
